# Update developer guide for integrated graph display
#
# 1) The cached "Fixed" date/time footer field (datetimeFigureOut) is
#    updated from "1/7/2017" to "11/1/17" wherever it is cached.
# 2) The "BrowserPanel" class box on the slide is renamed to
#    "GraphDisplay".

$p = $ppt.ActivePresentation

$oldDate = "1/7/2017"
$newDate = "11/1/17"

function Update-DatePlaceholder($container) {
    for ($i = 1; $i -le $container.Shapes.Count; $i++) {
        $shp = $container.Shapes.Item($i)
        if ($shp.HasTextFrame -eq -1 -and $shp.TextFrame.TextRange.Text -eq $oldDate) {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

$slide = $p.Slides.Item(1)

# --- Slide master's Date placeholder -------------------------------------
$master = $slide.Design.SlideMaster
Update-DatePlaceholder $master

# --- The slide layout actually used by the slide --------------------------
Update-DatePlaceholder $slide.CustomLayout

# --- Rename the BrowserPanel class box on the slide ------------------------
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTextFrame -eq -1 -and $shp.TextFrame.TextRange.Text -eq "BrowserPanel") {
        $shp.TextFrame.TextRange.Text = "GraphDisplay"
    }
}
